$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 115
$ws1.Range("F3").Value = 7588
$ws1.Range("F4").Value = 288
$ws1.Range("F5").Value = 23
$ws1.Range("F6").Value = 462
$ws1.Range("F7").Value = 4232
$ws1.Range("F9").Value = 589
$ws1.Range("F11").Value = 672

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 115
$ws4.Range("F4").Value = 7588
$ws4.Range("F6").Value = 288
$ws4.Range("F7").Value = 23
$ws4.Range("F8").Value = 462
$ws4.Range("F9").Value = 4232
$ws4.Range("F11").Value = 589
$ws4.Range("F13").Value = 672
